$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dates = @(
    "08-09-2021",
    "09-09-2021",
    "10-09-2021",
    "11-09-2021",
    "12-09-2021",
    "13-09-2021",
    "14-09-2021",
    "15-09-2021"
)

$startRow = 252
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    # These dates look like dd-mm-yyyy (day <= 15), so a plain Value
    # assignment gets silently parsed into a real date serial by Excel
    # instead of staying the literal text used throughout the rest of the
    # column. Force text interpretation with NumberFormat, then clear the
    # formatting back off so the cell ends up identical (no style index)
    # to the existing plain-text date cells above it.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
    $cell.ClearFormats()

    $ws.Cells.Item($row, 2).Value = 17537
    $ws.Cells.Item($row, 3).Value = 1456
    $ws.Cells.Item($row, 4).Value = 521
}
